# Insert a new price record for "Crimpson Seedless" grapes (Terminal
# Hortofrutícola Agro Chillán, weekly report) just above the existing
# row 160, shifting all the following rows (old 160..246) down by one
# row (to 161..247), and populate the newly opened row 160 with the
# new record's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 160:246 down to 161:247, opening up a blank row 160.
$ws.Rows(160).Insert()

# Fill in the new row 160 with the new daily price record.
$ws.Range("A160").Value = 7
$ws.Range("B160").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C160").Value = "Ñuble"
$ws.Range("D160").Value = 45086
$ws.Range("E160").Value = 16
$ws.Range("F160").Value = "Fruta"
$ws.Range("G160").Value = 100109
$ws.Range("H160").Value = "Uva"
$ws.Range("I160").Value = 100109001
$ws.Range("J160").Value = "Uva"
$ws.Range("K160").Value = "Crimpson Seedless"
$ws.Range("L160").Value = "Especial"
$ws.Range("M160").Value = 60
$ws.Range("N160").Value = 14000
$ws.Range("O160").Value = 14000
$ws.Range("P160").Value = 14000
$ws.Range("Q160").Value = "$/bandeja 18 kilos"
$ws.Range("R160").Value = "Región de O'Higgins"
$ws.Range("S160").Value = 778
$ws.Range("T160").Value = 18
